{"js": "// Load all body paragraphs so we can find the ones we need by their text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Helper: Office.js insertOoxml requires a full \"flat OPC\" package (a\n// <pkg:package> wrapper), not a bare WordprocessingML fragment.\nfunction flatOpcBody(bodyXml) {\n  return `<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">\n<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData>\n</pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>${bodyXml}</w:body></w:document></pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n}\n\n// 1) Wrap the run in the \"Etsyes\" paragraph with the proofing-error\n//    markers (<w:proofErr w:type=\"spellStart\"/> ... <w:proofErr\n//    w:type=\"spellEnd\"/>) Word writes around a word flagged by the\n//    spell checker.\nconst etsyesPara = items.find(p => p.text === \"Etsyes\");\nif (etsyesPara) {\n  const range = etsyesPara.getRange();\n\n  const startOoxml = flatOpcBody('<w:p><w:proofErr w:type=\"spellStart\"/></w:p>');\n  range.insertOoxml(startOoxml, Word.InsertLocation.start);\n\n  const endOoxml = flatOpcBody('<w:p><w:proofErr w:type=\"spellEnd\"/></w:p>');\n  range.insertOoxml(endOoxml, Word.InsertLocation.end);\n}\n\n// 2) Remove the \"Tests\", \"Linhas linhas\" and \"Linhas]\" paragraphs entirely.\nconst textsToRemove = new Set([\"Tests\", \"Linhas linhas\", \"Linhas]\"]);\nfor (const p of items) {\n  if (textsToRemove.has(p.text)) {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# WordprocessingML namespace used by the raw XML fragments below.\n$wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n# 1) Wrap the run in the \"Etsyes\" paragraph with the proofing-error\n#    markers (<w:proofErr w:type=\"spellStart\"/> ... <w:proofErr\n#    w:type=\"spellEnd\"/>) Word writes around a word flagged by the\n#    spell checker.\n#    Range.InsertXML *replaces* the addressed range's contents, so we\n#    rebuild the whole paragraph (ProofErr + original run) in one shot.\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"Etsyes\") {\n        $p.Range.InsertXML(\"<w:p $wNs><w:proofErr w:type=`\"spellStart`\"/><w:r><w:t>Etsyes</w:t></w:r><w:proofErr w:type=`\"spellEnd`\"/></w:p>\")\n        break\n    }\n}\n\n# 2) Remove the \"Tests\", \"Linhas linhas\" and \"Linhas]\" paragraphs\n#    entirely (walk backwards so deleting doesn't shift the indices of\n#    paragraphs still to be visited).\n$toRemove = @(\"Tests\", \"Linhas linhas\", \"Linhas]\")\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($toRemove -contains $text) {\n        $p.Range.Delete()\n    }\n}\n"}
